$d = $word.ActiveDocument

$replacements = @(
    @("485÷7=69, 2", "542÷7=77, 3"),
    @("243÷2=121, 1", "726÷3=242, 0"),
    @("558÷5=111, 3", "854÷9=94, 8"),
    @("115÷4=28, 3", "696÷5=139, 1"),
    @("532÷6=88, 4", "951÷7=135, 6"),
    @("819÷2=409, 1", "637÷6=106, 1"),
    @("361÷3=120, 1", "788÷8=98, 4"),
    @("328÷7=46, 6", "422÷3=140, 2"),
    @("632÷3=210, 2", "809÷7=115, 4"),
    @("160÷6=26, 4", "848÷7=121, 1"),
    @("679÷9=75, 4", "728÷2=364, 0"),
    @("151÷2=75, 1", "840÷3=280, 0"),
    @("526÷4=131, 2", "946÷9=105, 1"),
    @("463÷6=77, 1", "505÷6=84, 1"),
    @("687÷3=229, 0", "141÷6=23, 3"),
    @("400÷3=133, 1", "931÷6=155, 1"),
    @("189÷8=23, 5", "708÷3=236, 0"),
    @("679÷8=84, 7", "803÷5=160, 3"),
    @("294÷2=147, 0", "163÷3=54, 1"),
    @("374÷4=93, 2", "773÷9=85, 8"),
    @("414÷6=69, 0", "724÷9=80, 4"),
    @("814÷7=116, 2", "689÷3=229, 2"),
    @("825÷8=103, 1", "815÷5=163, 0"),
    @("980÷6=163, 2", "977÷6=162, 5"),
    @("149÷3=49, 2", "132÷3=44, 0")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
